$p = $ppt.ActivePresentation

# =========================================================================
# Slide 20 ("Compound Statement"): fix the curly-quote typo in the first
# bullet. The text reads:
#   A sequence of zero or more statements enclosed in braces "{" and "}".
# but the quote mark right after "{" was (incorrectly) a *left* quote mark
# instead of a *right* quote mark, so the quoting around "{" didn't match
# the quoting already used around "}". Fixing that one character forces a
# run split at that point (and at the "braces " boundary right before it),
# since the corrected text no longer shares identical run content with its
# neighbors.
# =========================================================================
$slide20 = $p.Slides.Item(20)
$shape20 = $slide20.Shapes.Item(2)
$tr20 = $shape20.TextFrame.TextRange
$para20 = $tr20.Paragraphs(1, 1)

$paraText20 = $para20.Text
$bracesIdx = $paraText20.IndexOf("braces ")
$openBraceIdx = $paraText20.IndexOf("{")
$strayQuoteIdx = $openBraceIdx + 1

# Give "braces <left-quote>" its own run boundary.
$seg1 = $tr20.Characters($para20.Start + $bracesIdx, $openBraceIdx - $bracesIdx + 1)
$seg1.Text = $seg1.Text

# The stray left quote mark right after "{" should be a right quote mark.
$quoteChar = $tr20.Characters($para20.Start + $strayQuoteIdx, 1)
$quoteChar.Text = "”"

# Give the fixed right-quote plus the following space its own run boundary.
$seg2 = $tr20.Characters($para20.Start + $strayQuoteIdx, 2)
$seg2.Text = $seg2.Text

# =========================================================================
# Slide 27 ("Procedures"): the procedureCallStmt grammar rule was spread
# across two paragraphs:
#   procedureCallStmt = procId "(" [ actualParameters ] ")"
#                       ";" .
# Join the trailing  ";" .  onto the end of the first paragraph, delete the
# now-redundant continuation paragraph, and shrink the grammar block's font
# size from 18.5 pt down to 17.5 pt.
# =========================================================================
$slide27 = $p.Slides.Item(27)
$shape27 = $slide27.Shapes.Item(2)
$tr27 = $shape27.TextFrame.TextRange

# Locate the "procedureCallStmt = ..." paragraph and its continuation.
$paraCount27 = $tr27.Paragraphs().Count
$callStmtParaIndex = -1
for ($i = 1; $i -le $paraCount27; $i++) {
    if ($tr27.Paragraphs($i, 1).Text.StartsWith("procedureCallStmt")) {
        $callStmtParaIndex = $i
        break
    }
}
$continuationParaIndex = $callStmtParaIndex + 1

$callStmtPara = $tr27.Paragraphs($callStmtParaIndex, 1)
$visibleLen = $callStmtPara.Length - 1
$tailChars = $tr27.Characters($callStmtPara.Start + $visibleLen - 6, 6)
$tailChars.Text = $tailChars.Text + ' ";" .'

$continuationPara = $tr27.Paragraphs($continuationParaIndex, 1)
$continuationPara.Delete()

# Shrink the font size of the (now two-paragraph) grammar block to 17.5 pt.
$grammarPara1 = $tr27.Paragraphs($callStmtParaIndex, 1)
$grammarPara1.Font.Size = 17.5

$grammarPara2 = $tr27.Paragraphs($callStmtParaIndex + 1, 1)
$grammarPara2.Font.Size = 17.5
